$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-valued columns (D, E, G) keep their cells formatted as Text
# so values are written as strings rather than being auto-converted to numbers/percentages.
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = '279.08'
$ws.Range("E2").Value = '0.90%'
$ws.Range("G2").Value = '20'
$ws.Range("E3").Value = '0.23%'
$ws.Range("G3").Value = '20'
$ws.Range("D4").Value = '4.837'
$ws.Range("E4").Value = '0.63%'
$ws.Range("G4").Value = '20'
$ws.Range("D5").Value = '0.06382'
$ws.Range("E5").Value = '0.44%'
$ws.Range("G5").Value = '20'
$ws.Range("D6").Value = '7.038'
$ws.Range("E6").Value = '1.21%'
$ws.Range("G6").Value = '20'
$ws.Range("D7").Value = '1.315'
$ws.Range("E7").Value = '2.34%'
$ws.Range("G7").Value = '20'
$ws.Range("D8").Value = '0.8945'
$ws.Range("E8").Value = '1.97%'
$ws.Range("G8").Value = '20'
$ws.Range("D9").Value = '0.1527'
$ws.Range("E9").Value = '0.19%'
$ws.Range("G9").Value = '20'
$ws.Range("D10").Value = '0.05823'
$ws.Range("E10").Value = '14.61%'
$ws.Range("G10").Value = '20'
$ws.Range("D11").Value = '0.07508'
$ws.Range("E11").Value = '-0.13%'
$ws.Range("G11").Value = '20'
$ws.Range("D12").Value = '0.02924'
$ws.Range("E12").Value = '-1.60%'
$ws.Range("G12").Value = '20'
$ws.Range("D13").Value = '0.08989'
$ws.Range("E13").Value = '-0.27%'
$ws.Range("G13").Value = '20'
$ws.Range("D14").Value = '0.001567'
$ws.Range("E14").Value = '0.31%'
$ws.Range("G14").Value = '20'
$ws.Range("D15").Value = '0.0006421'
$ws.Range("E15").Value = '0.52%'
$ws.Range("G15").Value = '20'
$ws.Range("E16").Value = '2.31%'
$ws.Range("G16").Value = '20'
$ws.Range("D17").Value = '3.475'
$ws.Range("E17").Value = '0.72%'
$ws.Range("G17").Value = '20'
$ws.Range("D18").Value = '3.308'
$ws.Range("E18").Value = '-0.03%'
$ws.Range("G18").Value = '20'
$ws.Range("D19").Value = '2.228'
$ws.Range("E19").Value = '-1.95%'
$ws.Range("G19").Value = '20'
$ws.Range("E20").Value = '-0.82%'
$ws.Range("G20").Value = '20'
$ws.Range("D21").Value = '0.1350'
$ws.Range("E21").Value = '1.05%'
$ws.Range("G21").Value = '20'
$ws.Range("D22").Value = '3.909'
$ws.Range("E22").Value = '-1.20%'
$ws.Range("G22").Value = '20'
$ws.Range("D23").Value = '0.04403'
$ws.Range("E23").Value = '-0.06%'
$ws.Range("G23").Value = '20'
$ws.Range("E24").Value = '8.89%'
$ws.Range("G24").Value = '20'
$ws.Range("D25").Value = '0.001177'
$ws.Range("E25").Value = '0.56%'
$ws.Range("G25").Value = '20'
$ws.Range("D26").Value = '0.004279'
$ws.Range("E26").Value = '10.87%'
$ws.Range("G26").Value = '20'
$ws.Range("G27").Value = '20'
$ws.Range("D28").Value = '0.0001180'
$ws.Range("E28").Value = '-1.67%'
$ws.Range("G28").Value = '20'
$ws.Range("D29").Value = '0.0001654'
$ws.Range("E29").Value = '-14.59%'
$ws.Range("G29").Value = '20'
$ws.Range("G30").Value = '20'
$ws.Range("G31").Value = '20'
$ws.Range("G32").Value = '20'
$ws.Range("G33").Value = '20'
$ws.Range("G34").Value = '20'
$ws.Range("G35").Value = '20'
$ws.Range("G36").Value = '20'
$ws.Range("G37").Value = '20'
$ws.Range("G38").Value = '20'
$ws.Range("G39").Value = '20'
$ws.Range("D40").Value = '0.04065'
$ws.Range("E40").Value = '-1.83%'
$ws.Range("G40").Value = '20'
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").Value = '0.006709'
$ws.Range("E41").Value = '-1.51%'
$ws.Range("G41").Value = '20'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = '0.1395'
$ws.Range("E42").Value = '18.51%'
$ws.Range("G42").Value = '20'
$ws.Range("D43").Value = '0.002129'
$ws.Range("E43").Value = '5.44%'
$ws.Range("G43").Value = '20'
$ws.Range("D44").Value = '0.01098'
$ws.Range("E44").Value = '-1.73%'
$ws.Range("G44").Value = '20'
$ws.Range("D45").Value = '0.00005563'
$ws.Range("E45").Value = '7.46%'
$ws.Range("G45").Value = '20'
$ws.Range("D46").Value = '1.561'
$ws.Range("E46").Value = '5.01%'
$ws.Range("G46").Value = '20'
$ws.Range("E47").Value = '-19.56%'
$ws.Range("G47").Value = '20'
$ws.Range("G48").Value = '20'
$ws.Range("G49").Value = '20'
$ws.Range("G50").Value = '20'
$ws.Range("G51").Value = '20'
